$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing #Genes values (column C) for several rows ---
$ws.Range("C15").Value = 5
$ws.Range("C17").Value = 1
$ws.Range("C18").Value = 12
$ws.Range("C19").Value = 4
$ws.Range("C20").Value = 3

# --- Append a brand-new module row (row 24) ---
$ws.Range("A24").Value = 22
$ws.Range("B24").Value = "Steroid hormone synthesis"
$ws.Range("C24").Value = 3
$ws.Range("D24").Value = 177
